$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.314316511154175
$ws.Range("B1").Value = 3.034732818603516
$ws.Range("C1").Value = 2.679712295532227
$ws.Range("D1").Value = 2.42784857749939
$ws.Range("E1").Value = 1.724100589752197
